$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Correct a handful of already-existing July (month=7) daily totals ---
$ws.Cells.Item(6,2).Value  = 24163.66   # day 7
$ws.Cells.Item(14,2).Value = 12214.2    # day 17
$ws.Cells.Item(15,2).Value = 7189.85    # day 18
$ws.Cells.Item(21,2).Value = 16418.74   # day 28
$ws.Cells.Item(22,2).Value = 25043.04   # day 29
$ws.Cells.Item(23,2).Value = 56663.15   # day 30

# --- 2) Insert a new row for July 31st, pushing everything below down by one ---
$ws.Rows.Item(24).Insert()
$ws.Cells.Item(24,1).Value = 31
$ws.Cells.Item(24,2).Value = 10799.55
$ws.Cells.Item(24,3).Value = 7
$ws.Cells.Item(24,4).Value = 2025
$ws.Cells.Item(24,5).Value = "07/2025"

# --- 3) Remove the whole April (month=4) block, which (after the insert above)
#        now lives in rows 67 through 86 ---
$ws.Range("A67:E86").EntireRow.Delete()
